# Sync attendance_reports: reorder "Recorded By" (column G) values so that
# any "System" / "system" entry in a comma-separated list is moved to the
# front of the list (list is reversed), for every row with data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 7).EntireColumn.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) {
    $lastRow = 157
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $text = [string]$val
    if ($text.IndexOf(",") -lt 0) {
        continue
    }

    $rawParts = $text.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $reversed = @()
        for ($i = $parts.Length - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $newText = [string]::Join(", ", $reversed)
        $cell.Value2 = $newText
    }
}
